$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 ---
$ws.Range("A5").Value = 94998123
$ws.Range("C5").Value = "Godkänd baserat på observatörens uppgifter"

# I5 holds a numeric-looking value but must stay text, like the sheet's other
# "Antal" cells; a leading apostrophe forces text entry, then resetting the
# style back to Normal avoids leaving a quote-prefix style on the cell.
$ws.Range("I5").Value = "'3"
$ws.Range("I5").Style = "Normal"

$ws.Range("P5").Value = "Örbäck, Hörendesjön, Vstm"
$ws.Range("S5").Value = 1
$ws.Range("X5").Value = ""
$ws.Range("AC5").Value = ""
$ws.Range("AW5").Value = "Zsombor Károlyi"
$ws.Range("AY5").Value = ""

# --- Row 6 ---
$ws.Range("A6").Value = 94998259
$ws.Range("C6").Value = "Godkänd baserat på observatörens uppgifter"

$ws.Range("I6").Value = "'7"
$ws.Range("I6").Style = "Normal"

$ws.Range("P6").Value = "Örbäck, Hörendesjön, Vstm"
$ws.Range("S6").Value = 1
$ws.Range("X6").Value = ""
$ws.Range("AC6").Value = ""
$ws.Range("AW6").Value = "Zsombor Károlyi"
$ws.Range("AY6").Value = ""
